# Auto-generated edit script applying the Seraph_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 562.7646999999999
$ws.Range("I19").Value = 499.16666
$ws.Range("J19").Value = 715.4
$ws.Range("K19").Value = 499.16666
$ws.Range("L19").Value = 715.4
$ws.Range("M19").Value = -324.16666
$ws.Range("N19").Value = -1065.4

# Row 43
$ws.Range("H43").Value = 6670.857
$ws.Range("J43").Value = 4499.6665
$ws.Range("L43").Value = 4499.6665
$ws.Range("N43").Value = -4637.6665

# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 112
$ws.Range("H112").Value = 2999.7856
$ws.Range("J112").Value = 3184.3845
$ws.Range("L112").Value = 9553.1535
$ws.Range("N112").Value = -11769.1535

# Row 132
$ws.Range("H132").Value = 1058.0416
$ws.Range("I132").Value = 1099.1818
$ws.Range("K132").Value = 3297.5454
$ws.Range("M132").Value = -767.5454

# Row 138
$ws.Range("H138").Value = 4725.2793
$ws.Range("I138").Value = 3446.6667
$ws.Range("J138").Value = 5734.7104
$ws.Range("K138").Value = 10340.0001
$ws.Range("L138").Value = 17204.1312
$ws.Range("M138").Value = -5200.000100000001
$ws.Range("N138").Value = -27484.1312

$ws = $wb.Worksheets.Item("ARM")
# Row 36
$ws.Range("H36").Value = 5821.5
$ws.Range("I36").Value = 2985.8
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 2985.8
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -2639.8
$ws.Range("N36").Value = -20692

# Row 61
$ws.Range("H61").Value = 1905.5555
$ws.Range("I61").Value = 1905.5555
$ws.Range("K61").Value = 1905.5555
$ws.Range("M61").Value = -1693.5555

# Row 74
$ws.Range("H74").Value = 1190.75
$ws.Range("I74").Value = 1162.6364
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1162.6364
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -288.6364000000001
$ws.Range("N74").Value = -3248

# Row 77
$ws.Range("H77").Value = 1190.75
$ws.Range("I77").Value = 1162.6364
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 5813.182000000001
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -1445.182000000001
$ws.Range("N77").Value = -16236

# Row 110
$ws.Range("H110").Value = 13997.8
$ws.Range("I110").Value = 14999.75
$ws.Range("J110").Value = 9990
$ws.Range("K110").Value = 14999.75
$ws.Range("L110").Value = 9990
$ws.Range("M110").Value = -12954.75
$ws.Range("N110").Value = -14080

# Row 122
$ws.Range("H122").Value = 3914.5334
$ws.Range("I122").Value = 3768.5833
$ws.Range("J122").Value = 4498.3335
$ws.Range("K122").Value = 11305.7499
$ws.Range("L122").Value = 13495.0005
$ws.Range("M122").Value = -8855.749899999999
$ws.Range("N122").Value = -18395.0005

# Row 136
$ws.Range("H136").Value = 1905.5555
$ws.Range("I136").Value = 1905.5555
$ws.Range("K136").Value = 5716.666499999999
$ws.Range("M136").Value = -3166.666499999999

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4846.143
$ws.Range("J20").Value = 1209.5
$ws.Range("L20").Value = 1209.5
$ws.Range("N20").Value = -1703.5

# Row 88
$ws.Range("H88").Value = 50114.332
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 50114.332
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 50114.332
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -50926.332

# Row 91
$ws.Range("H91").Value = 50114.332
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 50114.332
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 50114.332
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -52922.332

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4815.4375
$ws.Range("I31").Value = 4302.846
$ws.Range("K31").Value = 4302.846
$ws.Range("M31").Value = -4007.846

# Row 34
$ws.Range("H34").Value = 4815.4375
$ws.Range("I34").Value = 4302.846
$ws.Range("K34").Value = 4302.846
$ws.Range("M34").Value = -4100.846

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Row 86
$ws.Range("H86").Value = 7500
$ws.Range("I86").Value = 3750
$ws.Range("K86").Value = 3750
$ws.Range("M86").Value = -2627

# Row 89
$ws.Range("H89").Value = 7500
$ws.Range("I89").Value = 3750
$ws.Range("K89").Value = 18750
$ws.Range("M89").Value = -13134

# Row 107
$ws.Range("H107").Value = 1004.3333
$ws.Range("J107").Value = 1149
$ws.Range("L107").Value = 1149
$ws.Range("N107").Value = -4989

$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 12042.647
$ws.Range("I56").Value = 12042.647
$ws.Range("K56").Value = 12042.647
$ws.Range("M56").Value = -11512.647

# Row 86
$ws.Range("H86").Value = 300
$ws.Range("I86").Value = 300
$ws.Range("K86").Value = 900
$ws.Range("M86").Value = 286

# Row 89
$ws.Range("H89").Value = 300
$ws.Range("I89").Value = 300
$ws.Range("K89").Value = 2700
$ws.Range("M89").Value = 3228

# Row 103
$ws.Range("H103").Value = 158.8
$ws.Range("I103").Value = 158.8
$ws.Range("K103").Value = 476.4
$ws.Range("M103").Value = 402.6

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5255.1763
$ws.Range("I70").Value = 4121.727
$ws.Range("K70").Value = 4121.727
$ws.Range("M70").Value = -3851.727

# Row 73
$ws.Range("H73").Value = 5255.1763
$ws.Range("I73").Value = 4121.727
$ws.Range("K73").Value = 4121.727
$ws.Range("M73").Value = -3185.727

# Row 80
$ws.Range("H80").Value = 18674.75
$ws.Range("I80").Value = 9333
$ws.Range("J80").Value = 24279.8
$ws.Range("K80").Value = 9333
$ws.Range("L80").Value = 24279.8
$ws.Range("M80").Value = -8335
$ws.Range("N80").Value = -26275.8

# Row 83
$ws.Range("H83").Value = 18674.75
$ws.Range("I83").Value = 9333
$ws.Range("J83").Value = 24279.8
$ws.Range("K83").Value = 46665
$ws.Range("L83").Value = 121399
$ws.Range("M83").Value = -41673
$ws.Range("N83").Value = -131383

# Row 93
$ws.Range("H93").Value = 56494.25
$ws.Range("J93").Value = 56494.25
$ws.Range("L93").Value = 56494.25
$ws.Range("N93").Value = -60238.25

# Row 122
$ws.Range("H122").Value = 128202.625
$ws.Range("I122").Value = 3108.8333
$ws.Range("J122").Value = 503484
$ws.Range("K122").Value = 9326.499899999999
$ws.Range("L122").Value = 1510452
$ws.Range("M122").Value = -6876.499899999999
$ws.Range("N122").Value = -1515352

# Row 123
$ws.Range("H123").Value = 51250.375
$ws.Range("J123").Value = 51250.375
$ws.Range("L123").Value = 51250.375
$ws.Range("N123").Value = -56150.375

# Row 132
$ws.Range("H132").Value = 2700.5
$ws.Range("I132").Value = 2400.6155
$ws.Range("K132").Value = 7201.8465
$ws.Range("M132").Value = -4671.8465

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3568.4285
$ws.Range("I132").Value = 3568.4285
$ws.Range("K132").Value = 10705.2855
$ws.Range("M132").Value = -8175.2855

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1072.1666
$ws.Range("J113").Value = 1232
$ws.Range("L113").Value = 3696
$ws.Range("N113").Value = -8036
